$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 132, shifting existing rows 132-166 down to 133-167
$ws.Rows("132").Insert()

# Make sure the date column keeps the same format as the rest of column D
$ws.Range("D132").NumberFormat = $ws.Range("D133").NumberFormat

# Populate the new row 132 with data
$ws.Range("A132").Value = 4
$ws.Range("B132").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C132").Value = 'Los Lagos'
$ws.Range("D132").Value = 44511
$ws.Range("E132").Value = 10
$ws.Range("F132").Value = 100112021
$ws.Range("G132").Value = 'Ají'
$ws.Range("H132").Value = 'Inferno'
$ws.Range("I132").Value = 'Primera'
$ws.Range("J132").Value = 60
$ws.Range("K132").Value = 30000
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = 30000
$ws.Range("N132").Value = '$/caja 12 kilos'
$ws.Range("O132").Value = 'Región de Arica y Parinacota'
$ws.Range("P132").Value = 2500
$ws.Range("Q132").Value = 12
$ws.Range("R132").Value = 'Hortaliza'
